$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-26 03:48:11"
$ws.Range("M2").Value = "2.7 °C 3:04 TU"
$ws.Range("E3").Value = "2026-02-26 03:48:14"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "54%"
$ws.Range("O3").Value = "1.5 °C"
$ws.Range("E4").Value = "2026-02-26 03:48:16"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "92%"
$ws.Range("N4").Value = "4.6 °C 3:06 TU"
$ws.Range("O4").Value = "7.7 °C"
$ws.Range("E5").Value = "2026-02-26 03:48:18"
$ws.Range("O5").Value = "3.4 °C"
$ws.Range("E6").Value = "2026-02-26 03:48:21"
$ws.Range("N6").Value = "8.3 °C 3:00 TU"
$ws.Range("O6").Value = "9.7 °C"
$ws.Range("E7").Value = "2026-02-26 03:48:23"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "90%"
$ws.Range("N7").Value = "11.3 °C 3:25 TU"
$ws.Range("O7").Value = "11.8 °C"
$ws.Range("E8").Value = "2026-02-26 03:48:25"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "96%"
$ws.Range("N8").Value = "8.8 °C 3:28 TU"
$ws.Range("E9").Value = "2026-02-26 03:48:27"
$ws.Range("N9").Value = "10.3 °C 3:20 TU"
$ws.Range("O9").Value = "10.9 °C"
$ws.Range("E10").Value = "2026-02-26 03:48:30"
$ws.Range("N10").Value = "3.6 °C 3:22 TU"
$ws.Range("O10").Value = "4.4 °C"
$ws.Range("E11").Value = "2026-02-26 03:48:32"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "93%"
$ws.Range("N11").Value = "1.5 °C 3:03 TU"
$ws.Range("O11").Value = "2.2 °C"
$ws.Range("E12").Value = "2026-02-26 03:48:34"
$ws.Range("O12").Value = "9.7 °C"
$ws.Range("E13").Value = "2026-02-26 03:48:37"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "93%"
$ws.Range("J13").Value = "1031.5 hPa"
$ws.Range("L13").Value = "9.0 km/h - 61º 3:27 TU"
$ws.Range("N13").Value = "-1.9 °C 3:28 TU"
$ws.Range("O13").Value = "-0.6 °C"
$ws.Range("E14").Value = "2026-02-26 03:48:39"
$ws.Range("L14").Value = "20.9 km/h - 326º 3:09 TU"
$ws.Range("N14").Value = "9.2 °C 3:25 TU"
$ws.Range("O14").Value = "10.1 °C"
$ws.Range("E15").Value = "2026-02-26 03:48:41"
$ws.Range("N15").Value = "9.2 °C 3:15 TU"
$ws.Range("O15").Value = "10.6 °C"
$ws.Range("E16").Value = "2026-02-26 03:48:43"
$ws.Range("E17").Value = "2026-02-26 03:48:46"
$ws.Range("E18").Value = "2026-02-26 03:48:48"
$ws.Range("L18").Value = "4.3 km/h - 333º 3:18 TU"
$ws.Range("O18").Value = "8.5 °C"
$ws.Range("E19").Value = "2026-02-26 03:48:50"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "69%"
$ws.Range("E20").Value = "2026-02-26 03:48:52"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "57%"
$ws.Range("N20").Value = "-1.7 °C 3:09 TU"
$ws.Range("O20").Value = "0.6 °C"
$ws.Range("E21").Value = "2026-02-26 03:48:55"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "83%"
$ws.Range("J21").Value = "1028.5 hPa"
$ws.Range("N21").Value = "3.2 °C 3:29 TU"
$ws.Range("O21").Value = "4.6 °C"
$ws.Range("E22").Value = "2026-02-26 03:48:57"
$ws.Range("L22").Value = "15.8 km/h - 354º 3:29 TU"
$ws.Range("N22").Value = "-1.1 °C 3:07 TU"
$ws.Range("O22").Value = "0.3 °C"
$ws.Range("E23").Value = "2026-02-26 03:49:00"
$ws.Range("N23").Value = "1.4 °C 3:26 TU"
$ws.Range("O23").Value = "2.5 °C"
$ws.Range("E24").Value = "2026-02-26 03:49:02"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "75%"
$ws.Range("J24").Value = "1025.7 hPa"
$ws.Range("N24").Value = "3.4 °C 3:29 TU"
$ws.Range("O24").Value = "8.0 °C"
$ws.Range("E25").Value = "2026-02-26 03:49:04"
$ws.Range("O25").Value = "2.9 °C"
$ws.Range("E26").Value = "2026-02-26 03:49:07"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "46%"
$ws.Range("J26").Value = "1024.9 hPa"
$ws.Range("L26").Value = "12.2 km/h - 26º 3:24 TU"
$ws.Range("O26").Value = "7.0 °C"
$ws.Range("E27").Value = "2026-02-26 03:49:09"
$ws.Range("N27").Value = "1.7 °C 3:29 TU"
$ws.Range("O27").Value = "2.5 °C"
$ws.Range("E28").Value = "2026-02-26 03:49:11"
$ws.Range("J28").Value = "1026.0 hPa"
$ws.Range("L28").Value = "9.0 km/h - 244º 3:12 TU"
$ws.Range("N28").Value = "7.6 °C 3:29 TU"
$ws.Range("O28").Value = "8.7 °C"
$ws.Range("E29").Value = "2026-02-26 03:49:14"
$ws.Range("O29").Value = "9.9 °C"
$ws.Range("E30").Value = "2026-02-26 03:49:16"
$ws.Range("N30").Value = "10.5 °C 3:15 TU"
$ws.Range("O30").Value = "10.9 °C"
$ws.Range("E31").Value = "2026-02-26 03:49:18"
$ws.Range("O31").Value = "10.6 °C"
$ws.Range("E32").Value = "2026-02-26 03:49:21"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "73%"
$ws.Range("N32").Value = "0.2 °C 3:25 TU"
$ws.Range("O32").Value = "1.7 °C"
$ws.Range("E33").Value = "2026-02-26 03:49:23"
$ws.Range("J33").Value = "1029.0 hPa"
$ws.Range("N33").Value = "1.4 °C 3:29 TU"
$ws.Range("O33").Value = "2.9 °C"
$ws.Range("E34").Value = "2026-02-26 03:49:25"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "53%"
$ws.Range("L34").Value = "25.6 km/h - 49º 3:24 TU"
$ws.Range("M34").Value = "5.9 °C 3:18 TU"
$ws.Range("O34").Value = "2.0 °C"
$ws.Range("E35").Value = "2026-02-26 03:49:28"
$ws.Range("N35").Value = "8.2 °C 3:22 TU"
$ws.Range("O35").Value = "9.2 °C"
$ws.Range("E36").Value = "2026-02-26 03:49:30"
$ws.Range("E37").Value = "2026-02-26 03:49:32"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "88%"
$ws.Range("E38").Value = "2026-02-26 03:49:35"
$ws.Range("N38").Value = "6.2 °C 3:29 TU"
$ws.Range("O38").Value = "8.3 °C"
$ws.Range("E39").Value = "2026-02-26 03:49:37"
$ws.Range("E40").Value = "2026-02-26 03:49:39"
$ws.Range("J40").Value = "1029.3 hPa"
$ws.Range("N40").Value = "1.4 °C 3:12 TU"
$ws.Range("O40").Value = "2.4 °C"
$ws.Range("E41").Value = "2026-02-26 03:49:41"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "99%"
$ws.Range("E42").Value = "2026-02-26 03:49:43"
$ws.Range("O42").Value = "8.5 °C"
$ws.Range("E43").Value = "2026-02-26 03:49:46"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "96%"
$ws.Range("N43").Value = "1.6 °C 3:21 TU"
$ws.Range("O43").Value = "3.1 °C"
$ws.Range("E44").Value = "2026-02-26 03:49:48"
$ws.Range("N44").Value = "-1.6 °C 3:08 TU"
$ws.Range("O44").Value = "0.1 °C"
$ws.Range("E45").Value = "2026-02-26 03:49:50"
$ws.Range("J45").Value = "1027.1 hPa"
$ws.Range("N45").Value = "4.6 °C 3:09 TU"
$ws.Range("O45").Value = "6.4 °C"
$ws.Range("E46").Value = "2026-02-26 03:49:53"
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "97%"
$ws.Range("N46").Value = "6.5 °C 3:14 TU"
$ws.Range("O46").Value = "8.0 °C"
